$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.261.38"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "3.170.11"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'580.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'151.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.81%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.168.86"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  +6.20%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.80%  "
$ws.Range("D14").Value = "'37.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.80%  "
$ws.Range("D15").Value = "3.689.72"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").Value = "65.350.48"
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").Value = "'7.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("D18").Value = "3.170.75"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "'512.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.26%  "
$ws.Range("D21").Value = "'14.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.32%  "
$ws.Range("D22").Value = "'0.726"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.40%  "
$ws.Range("D23").Value = "'15.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.04%  "
$ws.Range("D24").Value = "'7.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").Value = "'85.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.63%  "
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.19%  "
$ws.Range("D30").Value = "'2.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +17.13%  "
$ws.Range("D31").Value = "'27.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.29%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("D34").Value = "'6.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.70%  "
$ws.Range("D35").Value = "'6.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("D36").Value = "'55.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'0.0908"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.66%  "
$ws.Range("D38").Value = "'475.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.68%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0423"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.33%  "
$ws.Range("D41").Value = "'8.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("D42").Value = "3.070.69"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  +2.73%  "
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.05%  "
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").Value = "'29.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.45%  "
$ws.Range("D47").Value = "0.0₃0602"
$ws.Range("E47").Value = "  +17.49%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").Value = "'2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.91%  "
$ws.Range("D51").Value = "'121.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
